# Applies the "don't need biomark marker tags as 'test tags' in metadata
# anymore" edit:
#  - AntennaMetadata (sheet1) becomes the active tab again (was TestTags),
#    its A column is widened, and the selection moves to B7.
#  - TestTags (sheet3) gets a new "notes" column (C) with per-tag remarks,
#    the old "Biomark Marker Tag" rows/type are dropped in favor of the
#    simplified "?" / "Pit tag on stick" typing plus notes, four stale rows
#    are removed (17 -> 13 total rows), and the selection moves to B7 too.

$wb = $excel.ActiveWorkbook

$wsAntenna = $wb.Worksheets.Item("AntennaMetadata")
$wsTags    = $wb.Worksheets.Item("TestTags")

# --- TestTags: rewrite the data body -----------------------------------
# Drop the four obsolete 999000000007xxx rows entirely and reshuffle the
# remaining tags/notes (17 data rows -> 13, incl. header) by deleting the
# now-unused trailing rows after the rewrite below.
#
# Shared-string table note: new strings get interned in first-write order,
# so the "4 detections..." note is written before the "notes" header (and
# both before any "no detections" cell) to match the workbook's string
# table ordering.
$wsTags.Cells.Item(4, 3).Value = "4 detections on 2020-10-08 at CF"
$wsTags.Range("C1").Value = "notes"

$tagRows = @(
    @(900230000102751, "Pit tag on stick", ""),
    @(900226001581072, "Pit tag on stick", ""),
    @(900230000004000, "", "4 detections on 2020-10-08 at CF"),
    @(900230000087405, "", "no detections"),
    @(900230000087408, "?", "shows up starting 2023-8-29 through 11-28 on all different different antennas "),
    @(900226001546996, "?", "10/31 - 11/28"),
    @(900230000088083, "", "5 detections 2023-5-18 at hp4 and cf5, cf6"),
    @(900230000087402, "?", "no detections"),
    @(900230000087403, "?", "no detections"),
    @(900230000088082, "", "45 detections 2023-5-18 at red barn only "),
    @(900230000228791, "", "24 detections; at RB on 10/7/2020, 7/27/2023 at confluence, again at red barn on 8/1/2023"),
    @(900230000087401, "", "55 detections in 2023: 6/6, 6/23. 7/25 at RB, CF and HP")
)

$r = 2
foreach ($row in $tagRows) {
    $wsTags.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne "") {
        $wsTags.Cells.Item($r, 2).Value = $row[1]
    } else {
        $wsTags.Cells.Item($r, 2).ClearContents()
    }
    $wsTags.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# The sheet used to have 17 rows (1 header + 16 data); only 13 remain now
# (1 header + 12 data), so remove the four leftover rows from the bottom.
$wsTags.Rows.Item(17).Delete()
$wsTags.Rows.Item(16).Delete()
$wsTags.Rows.Item(15).Delete()
$wsTags.Rows.Item(14).Delete()

# New notes column needs room; column B (type) also widens to fit "?" etc.
$wsTags.Columns.Item(2).ColumnWidth = 27.5
$wsTags.Columns.Item(3).ColumnWidth = 76

# --- AntennaMetadata: widen the site-name column ------------------------
$wsAntenna.Columns.Item(1).ColumnWidth = 41.5

# --- Selections -----------------------------------------------------------
$wsTags.Range("B7").Select() | Out-Null
$wsAntenna.Range("B7").Select() | Out-Null

# --- Active tab moves back to AntennaMetadata (was TestTags) -------------
$wsAntenna.Activate() | Out-Null
